$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 253
$ws.Range("I2").Value = 275
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 275
$ws.Range("L2").Value = 220
$ws.Range("M2").Value = -162
$ws.Range("N2").Value = -446
$ws.Range("H12").Value = 181.66667
$ws.Range("I12").Value = 198.66667
$ws.Range("J12").Value = 164.66667
$ws.Range("K12").Value = 198.66667
$ws.Range("L12").Value = 164.66667
$ws.Range("M12").Value = -28.66667000000001
$ws.Range("N12").Value = -504.66667
$ws.Range("H70").Value = 3459.8
$ws.Range("J70").Value = 3449.75
$ws.Range("L70").Value = 10349.25
$ws.Range("N70").Value = -10889.25
$ws.Range("H73").Value = 3459.8
$ws.Range("J73").Value = 3449.75
$ws.Range("L73").Value = 10349.25
$ws.Range("N73").Value = -12221.25
$ws.Range("H86").Value = 56334.668
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 56334.668
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H93").Value = 19000
$ws.Range("J93").Value = 19000
$ws.Range("L93").Value = 19000
$ws.Range("N93").Value = -23992
$ws.Range("H98").Value = 15834.333
$ws.Range("I98").Value = 13250
$ws.Range("K98").Value = 13250
$ws.Range("M98").Value = -11752
$ws.Range("H111").Value = 916.5
$ws.Range("I111").Value = 799.8
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 2399.4
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = 667.6000000000004
$ws.Range("N111").Value = -10634
$ws.Range("H122").Value = 15834.333
$ws.Range("I122").Value = 13250
$ws.Range("K122").Value = 39750
$ws.Range("M122").Value = -37300
$ws.Range("H138").Value = 2479.1304
$ws.Range("I138").Value = 1213.1428
$ws.Range("K138").Value = 3639.4284
$ws.Range("M138").Value = 1500.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 299.72223
$ws.Range("I2").Value = 310
$ws.Range("K2").Value = 310
$ws.Range("M2").Value = -197
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 150
$ws.Range("J5").Value = 150
$ws.Range("L5").Value = 150
$ws.Range("N5").Value = -374
$ws.Range("H32").Value = 1699.3334
$ws.Range("I32").Value = 749.5
$ws.Range("K32").Value = 749.5
$ws.Range("M32").Value = -462.5
$ws.Range("H45").Value = 1413.2075
$ws.Range("I45").Value = 1401.9231
$ws.Range("K45").Value = 1401.9231
$ws.Range("M45").Value = -1024.9231
$ws.Range("H74").Value = 2389.125
$ws.Range("I74").Value = 2030.4286
$ws.Range("J74").Value = 4900
$ws.Range("K74").Value = 2030.4286
$ws.Range("L74").Value = 4900
$ws.Range("M74").Value = -1156.4286
$ws.Range("N74").Value = -6648
$ws.Range("H77").Value = 2389.125
$ws.Range("I77").Value = 2030.4286
$ws.Range("J77").Value = 4900
$ws.Range("K77").Value = 10152.143
$ws.Range("L77").Value = 24500
$ws.Range("M77").Value = -5784.143
$ws.Range("N77").Value = -33236
$ws.Range("H97").Value = 1328.6471
$ws.Range("I97").Value = 845.1539
$ws.Range("K97").Value = 845.1539
$ws.Range("M97").Value = -349.1539
$ws.Range("H116").Value = 299.72223
$ws.Range("I116").Value = 310
$ws.Range("K116").Value = 310
$ws.Range("M116").Value = 1984
$ws.Range("H122").Value = 8333.333000000001
$ws.Range("I122").Value = 8333.333000000001
$ws.Range("K122").Value = 24999.999
$ws.Range("M122").Value = -22549.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 299.72223
$ws.Range("I3").Value = 310
$ws.Range("K3").Value = 310
$ws.Range("M3").Value = -196
$ws.Range("H4").Value = 150
$ws.Range("J4").Value = 150
$ws.Range("L4").Value = 150
$ws.Range("N4").Value = -380
$ws.Range("H94").Value = 1920.3334
$ws.Range("I94").Value = 1504.9375
$ws.Range("J94").Value = 2751.125
$ws.Range("K94").Value = 1504.9375
$ws.Range("L94").Value = 2751.125
$ws.Range("M94").Value = -1053.9375
$ws.Range("N94").Value = -3653.125
$ws.Range("H105").Value = 1369.8
$ws.Range("I105").Value = 1462.25
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 1462.25
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 284.75
$ws.Range("N105").Value = -4494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 24.333334
$ws.Range("I7").Value = 26.5
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 26.5
$ws.Range("L7").Value = 20
$ws.Range("M7").Value = 86.5
$ws.Range("N7").Value = -246
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1574
$ws.Range("H17").Value = 3008
$ws.Range("I17").Value = 3008
$ws.Range("K17").Value = 3008
$ws.Range("M17").Value = -2834
$ws.Range("H22").Value = 774
$ws.Range("I22").Value = 774
$ws.Range("K22").Value = 774
$ws.Range("M22").Value = -424
$ws.Range("H69").Value = 17036.2
$ws.Range("I69").Value = 8295.25
$ws.Range("K69").Value = 8295.25
$ws.Range("M69").Value = -7546.25
$ws.Range("H72").Value = 17036.2
$ws.Range("I72").Value = 8295.25
$ws.Range("K72").Value = 24885.75
$ws.Range("M72").Value = -21141.75
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 209.66667
$ws.Range("I107").Value = 84.666664
$ws.Range("K107").Value = 253.999992
$ws.Range("M107").Value = 1666.000008
$ws.Range("H131").Value = 4889.3335
$ws.Range("I131").Value = 4890
$ws.Range("J131").Value = 4888
$ws.Range("K131").Value = 14670
$ws.Range("L131").Value = 14664
$ws.Range("M131").Value = -9630
$ws.Range("N131").Value = -24744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 699.8
$ws.Range("I2").Value = 906.7273
$ws.Range("J2").Value = 130.75
$ws.Range("K2").Value = 906.7273
$ws.Range("L2").Value = 130.75
$ws.Range("M2").Value = -793.7273
$ws.Range("N2").Value = -356.75
$ws.Range("H70").Value = 1533
$ws.Range("I70").Value = 1533
$ws.Range("K70").Value = 1533
$ws.Range("M70").Value = -1263
$ws.Range("H73").Value = 1533
$ws.Range("I73").Value = 1533
$ws.Range("K73").Value = 1533
$ws.Range("M73").Value = -597
$ws.Range("H102").Value = 67301.836
$ws.Range("I102").Value = 80642.2
$ws.Range("J102").Value = 600
$ws.Range("K102").Value = 80642.2
$ws.Range("L102").Value = 600
$ws.Range("M102").Value = -79020.2
$ws.Range("N102").Value = -3844
$ws.Range("H107").Value = 1456
$ws.Range("I107").Value = 860.125
$ws.Range("J107").Value = 2409.4
$ws.Range("K107").Value = 860.125
$ws.Range("L107").Value = 2409.4
$ws.Range("M107").Value = 1059.875
$ws.Range("N107").Value = -6249.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1490.7273
$ws.Range("I16").Value = 1388.4445
$ws.Range("J16").Value = 1951
$ws.Range("K16").Value = 1388.4445
$ws.Range("L16").Value = 1951
$ws.Range("M16").Value = -1218.4445
$ws.Range("N16").Value = -2291
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("H136").Value = 900
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = -150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 2691.5
$ws.Range("I81").Value = 2868.5
$ws.Range("J81").Value = 2337.5
$ws.Range("K81").Value = 5737
$ws.Range("L81").Value = 4675
$ws.Range("M81").Value = -4676
$ws.Range("N81").Value = -6797
$ws.Range("H84").Value = 2691.5
$ws.Range("I84").Value = 2868.5
$ws.Range("J84").Value = 2337.5
$ws.Range("K84").Value = 28685
$ws.Range("L84").Value = 23375
$ws.Range("M84").Value = -23381
$ws.Range("N84").Value = -33983
$ws.Range("H107").Value = 999.6667
$ws.Range("I107").Value = 999.6667
$ws.Range("K107").Value = 2999.0001
$ws.Range("M107").Value = -1079.0001
$ws.Range("H109").Value = 69999
$ws.Range("J109").Value = 69999
$ws.Range("L109").Value = 69999
$ws.Range("N109").Value = -72773
$ws.Range("H122").Value = 3502
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H125").Value = 43500
$ws.Range("J125").Value = 43500
$ws.Range("L125").Value = 43500
$ws.Range("N125").Value = -53340
$ws.Range("H126").Value = 1033.3334
$ws.Range("I126").Value = 1033.3334
$ws.Range("K126").Value = 3100.0002
$ws.Range("M126").Value = -630.0001999999999
$ws.Range("H136").Value = 1155.125
$ws.Range("I136").Value = 891.7143
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 2675.1429
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -125.1428999999998
$ws.Range("N136").Value = -14097
